$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 139-148 (crossPlot.* keys rotated) ---
$ws.Range("A139").Value = "crossPlot.create-datapack"
$ws.Range("B139").Value = "Create Datapack"

$ws.Range("A140").Value = "crossPlot.time.select-unit"
$ws.Range("B140").Value = "Please select a unit."

$ws.Range("A141").Value = "crossPlot.time.no-chart-units-available"
$ws.Range("B141").Value = "No units available"

$ws.Range("A142").Value = "crossPlot.time.disabled-unit-reason"
$ws.Range("B142").Value = "You cannot change the age unit since it must be on the x-axis."

$ws.Range("A143").Value = "crossPlot.time.select-datapack"
$ws.Range("B143").Value = "Please select a datapack."

$ws.Range("A144").Value = "crossPlot.time.no-unit-selected"
$ws.Range("B144").Value = "No unit selected"

$ws.Range("A145").Value = "crossPlot.time.xAxis"
$ws.Range("B145").Value = "X Axis (Age-only)"

$ws.Range("A146").Value = "crossPlot.time.yAxis"
$ws.Range("B146").Value = "Y Axis (Age/Depth)"

$ws.Range("A147").Value = "crossPlot.sidebar.no-markers"
$ws.Range("B147").Value = " No Markers Available (Please Add a Marker by Clicking the Chart in the Main View)"

$ws.Range("A148").Value = "crossPlot.sidebar.no-models"
$ws.Range("B148").Value = "No Models Available (Please Add a Model by Clicking the Chart in the Main View)"

# --- Block 2: rows 273-296 (settings.datapacks.* keys rotated) ---
$ws.Range("A273").Value = "settings.datapacks.pdf-upload"
$ws.Range("B273").Value = "PDF Upload"
$ws.Range("C273").ClearContents()

$ws.Range("A274").Value = "settings.datapacks.seeMore"
$ws.Range("B274").Value = "See More..."
$ws.Range("C274").ClearContents()

$ws.Range("A275").Value = "settings.datapacks.seeLess"
$ws.Range("B275").Value = "See Less..."
$ws.Range("C275").ClearContents()

$ws.Range("A276").Value = "settings.datapacks.upload-form.title"
$ws.Range("B276").Value = "Upload Your Own Datapack"
$ws.Range("C276").Value = "上传你的数据包"

$ws.Range("A277").Value = "settings.datapacks.upload-form.no-file"
$ws.Range("B277").Value = "No file selected"
$ws.Range("C277").Value = "没有选择文件"

$ws.Range("A278").Value = "settings.datapacks.upload-form.name"
$ws.Range("B278").Value = "Datapack Name"
$ws.Range("C278").Value = "数据包名称"

$ws.Range("A279").Value = "settings.datapacks.upload-form.name-placeholder"
$ws.Range("B279").Value = "Enter a name for your datapack."
$ws.Range("C279").Value = "请为你的数据包输入名称"

$ws.Range("A280").Value = "settings.datapacks.upload-form.author"
$ws.Range("B280").Value = "Authored By"
$ws.Range("C280").Value = "作者"

$ws.Range("A281").Value = "settings.datapacks.upload-form.author-placeholder"
$ws.Range("B281").Value = "Credited to..."
$ws.Range("C281").Value = "作者为..."

$ws.Range("A282").Value = "settings.datapacks.upload-form.description"
$ws.Range("B282").Value = "Datapack Description"
$ws.Range("C282").Value = "数据包概述"

$ws.Range("A283").Value = "settings.datapacks.upload-form.description-placeholder"
$ws.Range("B283").Value = "Enter a description for your datapack."
$ws.Range("C283").Value = "请为你的数据包添加概述"

$ws.Range("A284").Value = "settings.datapacks.upload-form.tags"
$ws.Range("B284").Value = "Tags"
$ws.Range("C284").Value = "标签"

$ws.Range("A285").Value = "settings.datapacks.upload-form.make-public"
$ws.Range("B285").Value = "Make Datapack Publicly Accessible"
$ws.Range("C285").Value = "使数据包对公共可见"

$ws.Range("A286").Value = "settings.datapacks.upload-form.button.add-ref"
$ws.Range("B286").Value = "Add Reference"
$ws.Range("C286").Value = "添加引用"

$ws.Range("A287").Value = "settings.datapacks.upload-form.button.more"
$ws.Range("B287").Value = "More Options"
$ws.Range("C287").Value = "更多选项"

$ws.Range("A288").Value = "settings.datapacks.upload-form.button.finish"
$ws.Range("B288").Value = "Finish & Upload"
$ws.Range("C288").Value = "完成并上传"

$ws.Range("A289").Value = "settings.datapacks.upload-form.button.startover"
$ws.Range("B289").Value = "Start Over"
$ws.Range("C289").Value = "重新开始"

$ws.Range("A290").Value = "settings.datapacks.upload-form.reference"
$ws.Range("B290").Value = "Reference"
$ws.Range("C290").Value = "引用"

$ws.Range("A291").Value = "settings.datapacks.upload-form.contact"
$ws.Range("B291").Value = "Contact"
$ws.Range("C291").Value = "联系方式"

$ws.Range("A292").Value = "settings.datapacks.upload-form.contact-placeholder"
$ws.Range("B292").Value = "Enter your contact information"
$ws.Range("C292").Value = "输入你的联系方式"

$ws.Range("A293").Value = "settings.datapacks.upload-form.contact-helper-text"
$ws.Range("B293").Value = "(OPTIONAL) If you would like others to contact you about this datapack"
$ws.Range("C293").Value = "如果想要其他用户就此数据包问题联系你，请填写"

$ws.Range("A294").Value = "settings.datapacks.upload-form.notes"
$ws.Range("B294").Value = "Notes"
$ws.Range("C294").Value = "注释"

$ws.Range("A295").Value = "settings.datapacks.upload-form.notes-placeholder"
$ws.Range("B295").Value = "Enter notes for the datapack here"
$ws.Range("C295").Value = "为你的数据包添加注释"

$ws.Range("A296").Value = "settings.datapacks.upload-form.notes-helper-text"
$ws.Range("B296").Value = "(OPTIONAL) Generally notes are settings recommendations/How to use your datapack most efficiently"
$ws.Range("C296").Value = "如果想要指导其他用户如何有效地使用此数据包，请填写"
